# Auto-generated edit script applying scheduled market-price refresh
# to the Twintania_Profits workbook (8 crafting-job sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 17800
$ws.Range("J3").Value = 17800
$ws.Range("L3").Value = 17800
$ws.Range("N3").Value = -18028
$ws.Range("H13").Value = 9122.25
$ws.Range("J13").Value = 10496.667
$ws.Range("L13").Value = 10496.667
$ws.Range("N13").Value = -10834.667
$ws.Range("H76").Value = 4884.222
$ws.Range("I76").Value = 4807.25
$ws.Range("K76").Value = 4807.25
$ws.Range("M76").Value = -4492.25
$ws.Range("H79").Value = 4884.222
$ws.Range("I79").Value = 4807.25
$ws.Range("K79").Value = 4807.25
$ws.Range("M79").Value = -3715.25
$ws.Range("H100").Value = 501900.5
$ws.Range("I100").Value = 501900.5
$ws.Range("K100").Value = 501900.5
$ws.Range("M100").Value = -501359.5
$ws.Range("H102").Value = 17800
$ws.Range("J102").Value = 17800
$ws.Range("L102").Value = 17800
$ws.Range("N102").Value = -24290
$ws.Range("H127").Value = 3305
$ws.Range("I127").Value = 3505.7144
$ws.Range("J127").Value = 1900
$ws.Range("K127").Value = 10517.1432
$ws.Range("L127").Value = 5700
$ws.Range("M127").Value = -5557.143199999999
$ws.Range("N127").Value = -15620
$ws.Range("H129").Value = 3999.0625
$ws.Range("I129").Value = 1039.2
$ws.Range("J129").Value = 5344.4546
$ws.Range("K129").Value = 3117.6
$ws.Range("L129").Value = 16033.3638
$ws.Range("M129").Value = 1882.4
$ws.Range("N129").Value = -26033.3638
$ws.Range("H132").Value = 3348.5881
$ws.Range("I132").Value = 3005.1428
$ws.Range("J132").Value = 4951.3335
$ws.Range("K132").Value = 9015.428400000001
$ws.Range("L132").Value = 14854.0005
$ws.Range("M132").Value = -6485.428400000001
$ws.Range("N132").Value = -19914.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 13504.2
$ws.Range("I45").Value = 27449
$ws.Range("J45").Value = 4207.6665
$ws.Range("K45").Value = 27449
$ws.Range("L45").Value = 4207.6665
$ws.Range("M45").Value = -27072
$ws.Range("N45").Value = -4961.6665
$ws.Range("H132").Value = 7672.72
$ws.Range("I132").Value = 3825.75
$ws.Range("K132").Value = 11477.25
$ws.Range("M132").Value = -8947.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 36035.875
$ws.Range("J44").Value = 36035.875
$ws.Range("L44").Value = 36035.875
$ws.Range("N44").Value = -37029.875
$ws.Range("H105").Value = 3125.6538
$ws.Range("I105").Value = 3349.3
$ws.Range("J105").Value = 2380.1667
$ws.Range("K105").Value = 3349.3
$ws.Range("L105").Value = 2380.1667
$ws.Range("M105").Value = -1602.3
$ws.Range("N105").Value = -5874.1667
$ws.Range("H134").Value = 8804.134
$ws.Range("I134").Value = 5289.12
$ws.Range("K134").Value = 15867.36
$ws.Range("M134").Value = -13332.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1252.6875
$ws.Range("I31").Value = 968.0952
$ws.Range("J31").Value = 1796
$ws.Range("K31").Value = 968.0952
$ws.Range("L31").Value = 1796
$ws.Range("M31").Value = -673.0952
$ws.Range("N31").Value = -2386
$ws.Range("H34").Value = 1252.6875
$ws.Range("I34").Value = 968.0952
$ws.Range("J34").Value = 1796
$ws.Range("K34").Value = 968.0952
$ws.Range("L34").Value = 1796
$ws.Range("M34").Value = -766.0952
$ws.Range("N34").Value = -2200
$ws.Range("H99").Value = 11236.484
$ws.Range("I99").Value = 8815.333000000001
$ws.Range("J99").Value = 11774.519
$ws.Range("K99").Value = 8815.333000000001
$ws.Range("L99").Value = 11774.519
$ws.Range("M99").Value = -7317.333000000001
$ws.Range("N99").Value = -14770.519
$ws.Range("H126").Value = 11236.484
$ws.Range("I126").Value = 8815.333000000001
$ws.Range("J126").Value = 11774.519
$ws.Range("K126").Value = 26445.999
$ws.Range("L126").Value = 35323.557
$ws.Range("M126").Value = -23975.999
$ws.Range("N126").Value = -40263.557
$ws.Range("H132").Value = 31072.5
$ws.Range("I132").Value = 26481.682
$ws.Range("J132").Value = 35280.75
$ws.Range("K132").Value = 79445.046
$ws.Range("L132").Value = 105842.25
$ws.Range("M132").Value = -76915.046
$ws.Range("N132").Value = -110902.25
$ws.Range("H134").Value = 10483.5
$ws.Range("I134").Value = 8686.071
$ws.Range("K134").Value = 26058.213
$ws.Range("M134").Value = -23523.213

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 5999.6665
$ws.Range("I126").Value = 5999.6665
$ws.Range("K126").Value = 17998.9995
$ws.Range("M126").Value = -13058.9995
$ws.Range("H129").Value = 1029.6154
$ws.Range("I129").Value = 990.4167
$ws.Range("K129").Value = 2971.2501
$ws.Range("M129").Value = 2028.7499
$ws.Range("H130").Value = 2425
$ws.Range("I130").Value = 2425
$ws.Range("K130").Value = 7275
$ws.Range("M130").Value = -2255
$ws.Range("H131").Value = 38522.703
$ws.Range("I131").Value = 111919.22
$ws.Range("J131").Value = 1824.4445
$ws.Range("K131").Value = 335757.66
$ws.Range("L131").Value = 5473.333500000001
$ws.Range("M131").Value = -330717.66
$ws.Range("N131").Value = -15553.3335
$ws.Range("H134").Value = 1876.5
$ws.Range("I134").Value = 1876.5
$ws.Range("K134").Value = 5629.5
$ws.Range("M134").Value = -559.5
$ws.Range("H137").Value = 4050.5557
$ws.Range("I137").Value = 3491.6
$ws.Range("K137").Value = 10474.8
$ws.Range("M137").Value = -5374.799999999999
$ws.Range("H140").Value = 1267
$ws.Range("I140").Value = 1267
$ws.Range("K140").Value = 3801
$ws.Range("M140").Value = 1379
$ws.Range("H141").Value = 3932.25
$ws.Range("I141").Value = 3932.25
$ws.Range("K141").Value = 11796.75
$ws.Range("M141").Value = -6616.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 69.90909000000001
$ws.Range("I2").Value = 40.47059
$ws.Range("K2").Value = 40.47059
$ws.Range("M2").Value = 72.52941
$ws.Range("H99").Value = 17241.8
$ws.Range("I99").Value = 4636.857
$ws.Range("J99").Value = 46653.332
$ws.Range("K99").Value = 4636.857
$ws.Range("L99").Value = 46653.332
$ws.Range("M99").Value = -2390.857
$ws.Range("N99").Value = -51145.332
$ws.Range("H102").Value = 2944.25
$ws.Range("I102").Value = 2936.2856
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2936.2856
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1314.2856
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 4304
$ws.Range("I132").Value = 3852.125
$ws.Range("J132").Value = 5509
$ws.Range("K132").Value = 11556.375
$ws.Range("L132").Value = 16527
$ws.Range("M132").Value = -9026.375
$ws.Range("N132").Value = -21587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1160.4255
$ws.Range("I46").Value = 994.0714
$ws.Range("J46").Value = 1231
$ws.Range("K46").Value = 994.0714
$ws.Range("L46").Value = 1231
$ws.Range("M46").Value = -806.0714
$ws.Range("N46").Value = -1607
$ws.Range("H55").Value = 57.454544
$ws.Range("I55").Value = 60.2
$ws.Range("K55").Value = 60.2
$ws.Range("M55").Value = 112.8
$ws.Range("H122").Value = 3629.8572
$ws.Range("I122").Value = 3668
$ws.Range("K122").Value = 11004
$ws.Range("M122").Value = -8554
$ws.Range("H132").Value = 2793.8948
$ws.Range("I132").Value = 2436.8445
$ws.Range("J132").Value = 4132.8335
$ws.Range("K132").Value = 7310.5335
$ws.Range("L132").Value = 12398.5005
$ws.Range("M132").Value = -4780.5335
$ws.Range("N132").Value = -17458.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 23992.334
$ws.Range("J50").Value = 23992.334
$ws.Range("L50").Value = 23992.334
$ws.Range("N50").Value = -25254.334
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H61").Value = 40581.332
$ws.Range("I61").Value = 39725.273
$ws.Range("J61").Value = 49998
$ws.Range("K61").Value = 39725.273
$ws.Range("L61").Value = 49998
$ws.Range("M61").Value = -39433.273
$ws.Range("N61").Value = -50582
$ws.Range("H122").Value = 7766.3335
$ws.Range("I122").Value = 4577.4443
$ws.Range("K122").Value = 13732.3329
$ws.Range("M122").Value = -11282.3329

